# Auto-generated cell updates applying the scheduled market-data refresh diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I11").Value = 713.3333
$ws.Range("M11").Value = -573.3333
$ws.Range("H11").Value = 713.3333
$ws.Range("K11").Value = 713.3333
$ws.Range("I38").Value = 7554.533
$ws.Range("M38").Value = -22291.599
$ws.Range("H38").Value = 7707.375
$ws.Range("N38").Value = -30744
$ws.Range("K38").Value = 22663.599
$ws.Range("L38").Value = 30000
$ws.Range("J38").Value = 10000
$ws.Range("K46").Value = 36000
$ws.Range("I46").Value = 12000
$ws.Range("M46").Value = -35881
$ws.Range("H46").Value = 12000
$ws.Range("I48").Value = 11999
$ws.Range("M48").Value = -35705
$ws.Range("H48").Value = 13149.5
$ws.Range("N48").Value = -43484
$ws.Range("K48").Value = 35997
$ws.Range("L48").Value = 42900
$ws.Range("J48").Value = 14300
$ws.Range("L56").Value = 42900
$ws.Range("J56").Value = 14300
$ws.Range("I56").Value = 11999
$ws.Range("N56").Value = -43968
$ws.Range("M56").Value = -35463
$ws.Range("H56").Value = 13149.5
$ws.Range("K56").Value = 35997
$ws.Range("I60").Value = 12000
$ws.Range("M60").Value = -35516
$ws.Range("H60").Value = 12000
$ws.Range("K60").Value = 36000
$ws.Range("M74").Value = -4063
$ws.Range("H74").Value = 4999.5
$ws.Range("K74").Value = 4999
$ws.Range("N74").Value = -6872
$ws.Range("L74").Value = 5000
$ws.Range("J74").Value = 5000
$ws.Range("I74").Value = 4999
$ws.Range("I77").Value = 4999
$ws.Range("N77").Value = -34360
$ws.Range("M77").Value = -20315
$ws.Range("H77").Value = 4999.5
$ws.Range("K77").Value = 24995
$ws.Range("L77").Value = 25000
$ws.Range("J77").Value = 5000
$ws.Range("I82").Value = 3750
$ws.Range("M82").Value = -10844
$ws.Range("H82").Value = 3750
$ws.Range("K82").Value = 11250
$ws.Range("K85").Value = 11250
$ws.Range("I85").Value = 3750
$ws.Range("M85").Value = -9846
$ws.Range("H85").Value = 3750
$ws.Range("H113").Value = 1866.6666
$ws.Range("H135").Value = 2511.25
$ws.Range("N135").Value = -33613.5
$ws.Range("L135").Value = 28543.5
$ws.Range("J135").Value = 3171.5
$ws.Range("I138").Value = 992
$ws.Range("M138").Value = 2164
$ws.Range("H138").Value = 2687.2307
$ws.Range("N138").Value = -19266.3638
$ws.Range("K138").Value = 2976
$ws.Range("L138").Value = 8986.363799999999
$ws.Range("J138").Value = 2995.4546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("L2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("M2").Value = ""
$ws.Range("H2").Value = 0
$ws.Range("N2").Value = ""
$ws.Range("K2").Value = 0
$ws.Range("I41").Value = 15228.333
$ws.Range("M41").Value = -14814.333
$ws.Range("H41").Value = 15228.333
$ws.Range("K41").Value = 15228.333
$ws.Range("M74").Value = ""
$ws.Range("H74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("M77").Value = ""
$ws.Range("H77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M116").Value = ""
$ws.Range("H116").Value = 0
$ws.Range("N116").Value = ""
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("H132").Value = 5333.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I3").Value = 0
$ws.Range("M3").Value = ""
$ws.Range("H3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("N3").Value = ""
$ws.Range("L3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("I102").Value = 8750
$ws.Range("M102").Value = -5505
$ws.Range("H102").Value = 8750
$ws.Range("K102").Value = 8750
$ws.Range("H105").Value = 0
$ws.Range("N105").Value = ""
$ws.Range("L105").Value = 0
$ws.Range("J105").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 184
$ws.Range("N7").Value = -538.5
$ws.Range("K7").Value = 119.75
$ws.Range("L7").Value = 312.5
$ws.Range("J7").Value = 312.5
$ws.Range("I7").Value = 119.75
$ws.Range("M7").Value = -6.75
$ws.Range("I22").Value = 0
$ws.Range("M22").Value = ""
$ws.Range("H22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L31").Value = 2900
$ws.Range("J31").Value = 2900
$ws.Range("I31").Value = 2764.1428
$ws.Range("M31").Value = -2469.1428
$ws.Range("H31").Value = 2781.125
$ws.Range("N31").Value = -3490
$ws.Range("K31").Value = 2764.1428
$ws.Range("L34").Value = 2900
$ws.Range("J34").Value = 2900
$ws.Range("K34").Value = 2764.1428
$ws.Range("I34").Value = 2764.1428
$ws.Range("M34").Value = -2562.1428
$ws.Range("H34").Value = 2781.125
$ws.Range("N34").Value = -3304
$ws.Range("K107").Value = 626.5
$ws.Range("I107").Value = 626.5
$ws.Range("M107").Value = 1293.5
$ws.Range("H107").Value = 850.75
$ws.Range("N132").Value = -30147.5
$ws.Range("M132").Value = ""
$ws.Range("H132").Value = 8362.5
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 25087.5
$ws.Range("J132").Value = 8362.5
$ws.Range("I132").Value = 0
$ws.Range("J134").Value = 4428.4287
$ws.Range("I134").Value = 1660.8334
$ws.Range("N134").Value = -18355.2861
$ws.Range("M134").Value = -2447.5002
$ws.Range("H134").Value = 3151.077
$ws.Range("K134").Value = 4982.5002
$ws.Range("L134").Value = 13285.2861

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("L18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("I18").Value = 30
$ws.Range("N18").Value = ""
$ws.Range("M18").Value = 79
$ws.Range("H18").Value = 30
$ws.Range("K18").Value = 90
$ws.Range("N132").Value = -15630.5
$ws.Range("H132").Value = 1116.3334
$ws.Range("L132").Value = 10570.5
$ws.Range("J132").Value = 1174.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("K46").Value = 1700
$ws.Range("I46").Value = 1700
$ws.Range("M46").Value = -1544
$ws.Range("H46").Value = 1700
$ws.Range("I122").Value = 0
$ws.Range("M122").Value = ""
$ws.Range("H122").Value = 0
$ws.Range("K122").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I22").Value = 0
$ws.Range("M22").Value = ""
$ws.Range("H22").Value = 0
$ws.Range("N22").Value = ""
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("N27").Value = ""
$ws.Range("I27").Value = 0
$ws.Range("M27").Value = ""
$ws.Range("H27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("K46").Value = 4999.6665
$ws.Range("N46").Value = -13709
$ws.Range("L46").Value = 13333
$ws.Range("J46").Value = 13333
$ws.Range("I46").Value = 4999.6665
$ws.Range("M46").Value = -4811.6665
$ws.Range("H46").Value = 9166.333000000001
$ws.Range("M55").Value = -183.6
$ws.Range("H55").Value = 336.14285
$ws.Range("K55").Value = 356.6
$ws.Range("I55").Value = 356.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N62").Value = -6248
$ws.Range("L62").Value = 5000
$ws.Range("J62").Value = 5000
$ws.Range("H62").Value = 5000
$ws.Range("H65").Value = 5000
$ws.Range("N65").Value = -31240
$ws.Range("L65").Value = 25000
$ws.Range("J65").Value = 5000
$ws.Range("K107").Value = 2085
$ws.Range("N107").Value = -10218.9999
$ws.Range("L107").Value = 6378.999899999999
$ws.Range("J107").Value = 2126.3333
$ws.Range("I107").Value = 695
$ws.Range("M107").Value = -165
$ws.Range("H107").Value = 1553.8
$ws.Range("H126").Value = 6399
$ws.Range("K126").Value = 14995.0005
$ws.Range("N126").Value = -30440
$ws.Range("L126").Value = 25500
$ws.Range("J126").Value = 8500
$ws.Range("I126").Value = 4998.3335
$ws.Range("M126").Value = -12525.0005
